$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Title / date
Replace-Text "2026-02-09 Monday" "2026-02-10 Tuesday"

# Row 1
Replace-Text "381×6=" "973×9="
Replace-Text "771×2=" "228×7="
Replace-Text "550×3=" "783×9="
Replace-Text "367×8=" "358×7="
Replace-Text "631×9=" "969×7="

# Row 5
Replace-Text "862×6=" "385×5="
Replace-Text "717×5=" "295×8="
Replace-Text "872×8=" "896×3="
Replace-Text "180×9=" "779×4="
Replace-Text "968×6=" "846×8="

# Row 9
Replace-Text "618×2=" "142×2="
Replace-Text "268×3=" "853×2="
Replace-Text "334×7=" "796×2="
Replace-Text "810×2=" "879×9="
Replace-Text "383×3=" "405×2="

# Row 13
Replace-Text "945×7=" "445×5="
Replace-Text "670×5=" "112×2="
Replace-Text "643×9=" "882×9="
Replace-Text "579×7=" "464×9="
Replace-Text "764×6=" "740×6="

# Last row (20th row): cell-by-cell positional remap -
# old order: 490×7=, 241×6=, 415×3=, 209×2=, 968×8=
# new order: 582×9=, 610×8=, 241×6=, 920×2=, 680×7=
Replace-Text "490×7=" "582×9="
Replace-Text "241×6=" "610×8="
Replace-Text "415×3=" "241×6="
Replace-Text "209×2=" "920×2="
Replace-Text "968×8=" "680×7="
